$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: replace with what used to be row 3's data, plus a new boolean in F1
$ws.Range("A1").Value = "c"
$ws.Range("B1").Value = "t1"
$ws.Range("C1").Value = "Televisor"
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 0
$ws.Range("F1").Value = $true

# Row 2: keep A2/B2/C2, update the numeric columns
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 58
$ws.Range("F2").Value = $true

# Row 3 is removed entirely (its old data was moved up into row 1)
$ws.Rows("3:3").Delete()
